# Sprint Backlog: add "Spring 9 (M9)" sheet at the end of the workbook and
# populate it with the M9 task list, mirroring the layout of the preceding
# "Sprint 8 (M8)" sheet. Also updates the previously-active sheet's selection
# now that it is no longer the active tab.

$wb = $excel.ActiveWorkbook

# --- Sprint 8 (M8): it is no longer the active/selected tab -------------
$ws8 = $wb.Worksheets.Item(8)
$ws8.Range("A1:G5").Select() | Out-Null

# --- Add the new sheet at the very end of the tab strip -----------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws9 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws9.Name = "Spring 9 (M9)"

# --- Header row -----------------------------------------------------------
$ws9.Range("A1").Value = "Tasks"
$ws9.Range("B1").Value = "Responsible"
$ws9.Range("C1").Value = "Status"
$ws9.Range("D1").Value = 1
$ws9.Range("E1").Value = 2
$ws9.Range("F1").Value = 3
$ws9.Range("G1").Value = 4
$ws9.Range("A1:G1").Font.Bold = $true
$ws9.Range("A1:G1").WrapText = $true
$ws9.Rows.Item(1).RowHeight = 30

# --- Task rows (write column B before A so new shared strings land in ----
# --- the same order Excel produced them: "Done" before "Method           --
# --- Contracts ...") -------------------------------------------------------
$ws9.Range("B2").Value = "Done"
$ws9.Range("A2").Value = "Method Contracts (individual)"

$ws9.Range("A3").Value = "Ensure shipyard available at tech planets/Implement Option to Upgrade Ship"
$ws9.Range("B3").Value = "Stephen"

$ws9.Range("A4").Value = "Enforce money/slot limits for upgrading ship. Enforce some gadgets only available on higher tech level planet"
$ws9.Range("B4").Value = "Naman"

$ws9.Range("A5").Value = "Apply design pattern to code"
$ws9.Range("B5").Value = "Bhavesh"

$ws9.Range("A6").Value = "Code Critique and Java Doc"
$ws9.Range("B6").Value = "Hunter"

# --- Column width / selection on the new sheet ----------------------------
$ws9.Columns.Item(1).ColumnWidth = 89.5
$ws9.Range("A8").Select() | Out-Null

$ws9.Activate() | Out-Null
